$wb = $excel.ActiveWorkbook

# The workbook has two sheets that duplicate the same exhibition data:
# "展览" (Exhibitions) and "全部类型" (All types). Both need their
# "想去人数" (want-to-go count) column F updated for rows 2 and 3.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 739
    $ws.Range("F3").Value = 4121
}
